$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point drift on A21 (same instant, re-serialized)
$ws.Range("A21").Value = 45877.87534339121

# Append new row 22 with the latest sensor reading
$ws.Range("A22").Value = 45877.91685341245
$ws.Range("B22").Value = 2025
$ws.Range("C22").Value = 32
$ws.Range("D22").Value = 14.33
$ws.Range("E22").Value = 90.41
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0.28
$ws.Range("H22").Value = "NNE"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "22:00:16"

# Match the date-format style used on A21
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat
